# Edit: reverse the order of the "Periodo Mora" / "Valor Mora" data rows
# (rows 16-22, columns E and F) on sheet "Hoja1".
#
# Before: E16:E22 = 2211,2212,2301,2302,2303,2304,2305
#         F16:F22 = 40000,40000,40000,40000,40000,40000,32000
# After:  E16:E22 = 2305,2304,2303,2302,2301,2212,2211
#         F16:F22 = 32000,40000,40000,40000,40000,40000,40000

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Hoja1")

$firstRow = 16
$lastRow = 22

# Capture the current values for the two columns that need to be reversed.
$periodos = @()
$valores = @()
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $periodos += $ws.Cells.Item($r, 5).Value2
    $valores += $ws.Cells.Item($r, 6).Value2
}

# Write the values back in reverse order.
$count = $lastRow - $firstRow + 1
for ($i = 0; $i -lt $count; $i++) {
    $targetRow = $firstRow + $i
    $sourceIndex = $count - 1 - $i
    $ws.Cells.Item($targetRow, 5).Value = $periodos[$sourceIndex]
    $ws.Cells.Item($targetRow, 6).Value = $valores[$sourceIndex]
}

$wb.Save()
